# Update COVID country data ("paises.xlsx") with the later refresh (10:22 -> 10:52).
# The table (rows 4-216) is kept sorted descending by column B ("Casos totales"),
# so refreshed case counts for a few countries pushed them above their former
# neighbours -> those rows shift down by one (country name + all stats), while
# the country that overtook its neighbours gets the new top-of-block figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 10:52"

# Austria (row 20) - values refreshed in place, no re-sort needed
$ws.Range("B20").Value = 14102
$ws.Range("C20").Value = 61
$ws.Range("D20").Value = 7633
$ws.Range("E20").Value = 6085
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 384

# Rumania (row 31) - values refreshed in place, no re-sort needed
$ws.Range("E31").Value = 5387
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 332

# Block around rows 47-50: Bielorrusia overtakes Catar / Republica Dominicana / Finlandia
$ws.Range("A47").Value = "Bielorrusia"
$ws.Range("B47").Value = 3281
$ws.Range("C47").Value = 362
$ws.Range("D47").Value = 203
$ws.Range("E47").Value = 3045
$ws.Range("F47").Value = 57
$ws.Range("G47").Value = 4
$ws.Range("H47").Value = 33

$ws.Range("A48").Value = "Catar"
$ws.Range("B48").Value = 3231
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 334
$ws.Range("E48").Value = 2890
$ws.Range("F48").Value = 37
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 7

$ws.Range("A49").Value = "Republica Dominicana"
$ws.Range("B49").Value = 3167
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 152
$ws.Range("E49").Value = 2838
$ws.Range("F49").Value = 147
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 177

$ws.Range("A50").Value = "Finlandia"
$ws.Range("B50").Value = 3064
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 300
$ws.Range("E50").Value = 2705
$ws.Range("F50").Value = 74
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 59

# Block around rows 76-81: Banglades overtakes Hong Kong / Rep. Macedonia / Camerun / Eslovaquia / Oman
$ws.Range("A76").Value = "Banglades"
$ws.Range("B76").Value = 1012
$ws.Range("C76").Value = 209
$ws.Range("D76").Value = 42
$ws.Range("E76").Value = 924
$ws.Range("F76").Value = 1
$ws.Range("G76").Value = 7
$ws.Range("H76").Value = 46

$ws.Range("A77").Value = "Hong Kong"
$ws.Range("B77").Value = 1010
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 397
$ws.Range("E77").Value = 609
$ws.Range("F77").Value = 13
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 4

$ws.Range("A78").Value = "Republica de Macedonia"
$ws.Range("B78").Value = 854
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 44
$ws.Range("E78").Value = 772
$ws.Range("F78").Value = 15
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 38

$ws.Range("A79").Value = "Camerun"
$ws.Range("B79").Value = 848
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 98
$ws.Range("E79").Value = 738
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 12

$ws.Range("A80").Value = "Eslovaquia"
$ws.Range("B80").Value = 816
$ws.Range("C80").Value = 47
$ws.Range("D80").Value = 107
$ws.Range("E80").Value = 707
$ws.Range("F80").Value = 5
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 2

$ws.Range("A81").Value = "Oman"
$ws.Range("B81").Value = 813
$ws.Range("C81").Value = 86
$ws.Range("D81").Value = 130
$ws.Range("E81").Value = 679
$ws.Range("F81").Value = 3
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 4

# Block around rows 84-86: Afganistan overtakes Crucero / Bulgaria
$ws.Range("A84").Value = "Afganistan"
$ws.Range("B84").Value = 714
$ws.Range("C84").Value = 49
$ws.Range("D84").Value = 40
$ws.Range("E84").Value = 651
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 23

$ws.Range("A85").Value = "Crucero"
$ws.Range("B85").Value = 712
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 639
$ws.Range("E85").Value = 61
$ws.Range("F85").Value = 7
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 12

$ws.Range("A86").Value = "Bulgaria"
$ws.Range("B86").Value = 695
$ws.Range("C86").Value = 10
$ws.Range("D86").Value = 81
$ws.Range("E86").Value = 579
$ws.Range("F86").Value = 36
$ws.Range("G86").Value = 3
$ws.Range("H86").Value = 35
